$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J, matching style of existing header cells (H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-16
$data = @{
    2  = @(6, 7)
    3  = @(8, 8)
    4  = @(9, 9)
    5  = @(8, 8)
    6  = @(5, 6)
    7  = @(6, 7)
    8  = @(7, 8)
    9  = @(5, 5)
    10 = @(4, 4)
    11 = @(6, 6)
    12 = @(5, 5)
    13 = @(8, 9)
    14 = @(6, 6)
    15 = @(7, 7)
    16 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
